$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.811.76'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '3.151.00'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'579.04"
$ws.Range('E5').Value = '  +1.34%  '
$ws.Range('D6').Value = "'149.02"
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.148.13'
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').Value = "'0.526"
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('D10').Value = "'0.158"
$ws.Range('E10').Value = '  -2.23%  '
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').Value = "'0.500"
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('E13').Value = '  +2.89%  '
$ws.Range('D14').Value = "'37.15"
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('D15').Value = '3.672.11'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('D16').Value = '64.882.45'
$ws.Range('D17').Value = '3.169.09'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = "'7.13"
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').Value = "'503.65"
$ws.Range('E20').Value = '  -2.73%  '
$ws.Range('D21').Value = "'14.88"
$ws.Range('E21').Value = '  -0.77%  '
$ws.Range('D22').Value = "'0.715"
$ws.Range('E22').Value = '  -2.89%  '
$ws.Range('D23').Value = "'15.17"
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').Value = "'7.72"
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('D25').Value = "'84.36"
$ws.Range('E25').Value = '  -1.09%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').Value = "'2.92"
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = "'8.92"
$ws.Range('E28').Value = '  +1.99%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  +4.09%  '
$ws.Range('D31').Value = "'27.59"
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').Value = "'1.19"
$ws.Range('E33').Value = '  +0.63%  '
$ws.Range('D34').Value = "'6.23"
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('D35').Value = "'6.47"
$ws.Range('E35').Value = '  -1.52%  '
$ws.Range('D36').Value = "'54.75"
$ws.Range('E36').Value = '  -1.78%  '
$ws.Range('D37').Value = "'484.26"
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('D38').Value = "'0.0890"
$ws.Range('E38').Value = '  +2.79%  '
$ws.Range('D39').Value = "'0.0416"
$ws.Range('E39').Value = '  -1.36%  '
$ws.Range('D40').Value = "'2.92"
$ws.Range('E40').Value = '  -1.43%  '
$ws.Range('D41').Value = "'8.74"
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('D42').Value = '2.993.13'
$ws.Range('E42').Value = '  -3.83%  '
$ws.Range('E43').Value = '  -3.62%  '
$ws.Range('D44').Value = "'2.43"
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').Value = "'0.282"
$ws.Range('E45').Value = '  -4.43%  '
$ws.Range('D46').Value = "'28.11"
$ws.Range('E46').Value = '  -3.83%  '
$ws.Range('D47').Value = '0.0₃0587'
$ws.Range('E47').Value = '  +1.81%  '
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('D50').Value = "'2.23"
$ws.Range('E50').Value = '  -2.81%  '
$ws.Range('D51').Value = "'2.47"
$ws.Range('E51').Value = '  +13.53%  '
